$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect the unified "DataNode" concept (was "Property1")
$ws.Name = "DataNode"

# Row-height touch-ups that came along with the re-save (header row now wraps to
# two lines, the long description row shrank slightly)
$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(8).RowHeight = 81

# Selection left on H13 after editing
$ws.Range("H13").Select() | Out-Null
